$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number (45181 -> 2023-09-12)
# for every data row (2 through 357). Bump it by one day (45182 -> 2023-09-13)
# for all of them, matching the source diff exactly.
$range = $ws.Range("C2:C357")
for ($i = 1; $i -le $range.Rows.Count; $i++) {
    $cell = $range.Cells.Item($i, 1)
    if ($cell.Value2 -eq 45181) {
        $cell.Value2 = 45182
    }
}
